# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly scraped numbers, per the commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) -> row number : new F value
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 82
    6  = 39
    7  = 28
    9  = 7
    10 = 14406
    11 = 152
    13 = 5760
    14 = 586
    16 = 42
    18 = 1234
    21 = 177
    23 = 2930
    25 = 10546
    30 = 242
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "全部类型" (all types) -> row number : new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 82
    5  = 3017
    7  = 39
    8  = 28
    10 = 7
    11 = 14406
    12 = 152
    14 = 5760
    15 = 586
    17 = 42
    19 = 1234
    22 = 177
    24 = 2930
    27 = 10546
    32 = 242
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
